$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(3, 6).Value = 1287   # F3: 1282 -> 1287
$ws.Cells.Item(4, 6).Value = 88   # F4: 87 -> 88
$ws.Cells.Item(6, 6).Value = 124   # F6: 123 -> 124
$ws.Cells.Item(7, 6).Value = 851   # F7: 847 -> 851
$ws.Cells.Item(11, 6).Value = 1070   # F11: 1069 -> 1070
$ws.Cells.Item(12, 6).Value = 807   # F12: 806 -> 807
$ws.Cells.Item(14, 6).Value = 681   # F14: 677 -> 681
$ws.Cells.Item(15, 6).Value = 1322   # F15: 1309 -> 1322
$ws.Cells.Item(16, 6).Value = 1025   # F16: 1026 -> 1025
$ws.Cells.Item(19, 6).Value = 731   # F19: 730 -> 731
$ws.Cells.Item(23, 6).Value = 636   # F23: 635 -> 636
$ws.Cells.Item(24, 6).Value = 1221   # F24: 1220 -> 1221
$ws.Cells.Item(25, 6).Value = 142   # F25: 140 -> 142
$ws.Cells.Item(28, 6).Value = 5151   # F28: 5136 -> 5151
$ws.Cells.Item(29, 6).Value = 249   # F29: 245 -> 249
$ws.Cells.Item(31, 6).Value = 2421   # F31: 2419 -> 2421
$ws.Cells.Item(32, 6).Value = 5808   # F32: 5801 -> 5808
$ws.Cells.Item(34, 6).Value = 973   # F34: 970 -> 973
$ws.Cells.Item(35, 6).Value = 590   # F35: 588 -> 590
$ws.Cells.Item(37, 6).Value = 1041   # F37: 1040 -> 1041
$ws.Cells.Item(41, 6).Value = 667   # F41: 663 -> 667
$ws.Cells.Item(43, 6).Value = 30   # F43: 29 -> 30

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(8, 6).Value = 113   # F8: 112 -> 113
$ws.Cells.Item(12, 6).Value = 90   # F12: 91 -> 90
$ws.Cells.Item(15, 6).Value = 659   # F15: 658 -> 659
$ws.Cells.Item(32, 6).Value = 146   # F32: 144 -> 146
$ws.Cells.Item(33, 6).Value = 98   # F33: 97 -> 98
$ws.Cells.Item(38, 6).Value = 14   # F38: 13 -> 14
$ws.Cells.Item(41, 6).Value = 485   # F41: 482 -> 485

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(5, 6).Value = 747   # F5: 746 -> 747
$ws.Cells.Item(6, 6).Value = 364   # F6: 361 -> 364
$ws.Cells.Item(7, 6).Value = 209   # F7: 208 -> 209

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(5, 6).Value = 1287   # F5: 1282 -> 1287
$ws.Cells.Item(6, 6).Value = 88   # F6: 87 -> 88
$ws.Cells.Item(7, 6).Value = 364   # F7: 361 -> 364
$ws.Cells.Item(8, 6).Value = 209   # F8: 208 -> 209
$ws.Cells.Item(9, 6).Value = 209   # F9: 208 -> 209
$ws.Cells.Item(12, 6).Value = 851   # F12: 847 -> 851
$ws.Cells.Item(13, 6).Value = 113   # F13: 112 -> 113
$ws.Cells.Item(17, 6).Value = 1071   # F17: 1069 -> 1071
$ws.Cells.Item(18, 6).Value = 807   # F18: 806 -> 807
$ws.Cells.Item(21, 6).Value = 681   # F21: 677 -> 681
$ws.Cells.Item(22, 6).Value = 1322   # F22: 1309 -> 1322
$ws.Cells.Item(24, 6).Value = 1025   # F24: 1026 -> 1025
$ws.Cells.Item(27, 6).Value = 731   # F27: 730 -> 731
$ws.Cells.Item(30, 6).Value = 636   # F30: 635 -> 636
$ws.Cells.Item(31, 6).Value = 1221   # F31: 1220 -> 1221
$ws.Cells.Item(32, 6).Value = 142   # F32: 140 -> 142
$ws.Cells.Item(35, 6).Value = 5152   # F35: 5136 -> 5152
$ws.Cells.Item(36, 6).Value = 249   # F36: 245 -> 249
$ws.Cells.Item(38, 6).Value = 2421   # F38: 2419 -> 2421
$ws.Cells.Item(39, 6).Value = 5808   # F39: 5801 -> 5808
$ws.Cells.Item(40, 6).Value = 973   # F40: 970 -> 973
$ws.Cells.Item(42, 6).Value = 590   # F42: 588 -> 590
$ws.Cells.Item(44, 6).Value = 1041   # F44: 1040 -> 1041
$ws.Cells.Item(46, 6).Value = 667   # F46: 663 -> 667
$ws.Cells.Item(48, 6).Value = 30   # F48: 29 -> 30
$ws.Cells.Item(50, 6).Value = 485   # F50: 482 -> 485
